$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    5  = @{ E = 173; F = 116; H = 127 }
    10 = @{ E = 742; F = 411; H = 506 }
    11 = @{ E = 495; F = 276; H = 341 }
    12 = @{ E = 765; F = 452; H = 538 }
    14 = @{ E = 150; F = 84;  H = 118 }
    15 = @{ E = 212 }
    16 = @{ E = 240; F = 137; H = 185 }
    17 = @{ E = 128; F = 70;  H = 94 }
    18 = @{ E = 65 }
    20 = @{ E = 103; F = 48;  H = 85 }
    24 = @{ E = 288 }
    25 = @{ E = 350; F = 189; H = 249 }
    26 = @{ E = 225; F = 130; H = 155 }
    27 = @{ E = 405; F = 220; H = 302 }
    28 = @{ F = 114; H = 166 }
    29 = @{ E = 201 }
    30 = @{ E = 270; F = 167; H = 220 }
    33 = @{ E = 348; F = 184; H = 275 }
    34 = @{ E = 267 }
    36 = @{ E = 92 }
    37 = @{ E = 204 }
    39 = @{ E = 209 }
    40 = @{ E = 319 }
    41 = @{ E = 456; F = 228; H = 320 }
    42 = @{ E = 487; F = 271; H = 332 }
    43 = @{ E = 150 }
    44 = @{ E = 390; F = 206; H = 274 }
    45 = @{ E = 190 }
    46 = @{ E = 409; F = 238; H = 302 }
    47 = @{ E = 566; F = 314; H = 406 }
    48 = @{ E = 286; F = 137; H = 181 }
    49 = @{ E = 347; F = 171; H = 258 }
    50 = @{ E = 290 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
